# Scene 27C.docx edit
#
# The canonical-XML diff for this commit shows word/styles.xml gaining a
# second <w:style styleId="Subtitle"> definition: an extra "Subtitle"
# paragraph style (italic Georgia, dark-grey, 24pt, spaced like the other
# heading-ish styles) appended right after the "Subtitle" style that is
# already in the document. (The other style-boilerplate the diff shows --
# Heading1-6/Title/the stub Normal&TableNormal entries -- is already
# present in this copy of the document.)
#
# Word's Styles collection can't hold two styles that share a StyleId, so
# we recreate the duplicate the way the UI would if you dragged/copied a
# second "Subtitle" style in from another document: a new style entry
# with its own internal id, whose visible name is set to "Subtitle" to
# match, carrying identical formatting to the original.

$d = $word.ActiveDocument

$dup = $d.Styles.Add("Subtitle1", 1)
$dup.BaseStyle = $d.Styles.Item("Normal")
$dup.NextParagraphStyle = $d.Styles.Item("Normal")

$dup.ParagraphFormat.KeepWithNext = $true
$dup.ParagraphFormat.KeepTogether = $true
$dup.ParagraphFormat.SpaceBefore = 18
$dup.ParagraphFormat.SpaceAfter = 4
$dup.ParagraphFormat.LineSpacingRule = 5

$dup.Font.NameAscii = "Georgia"
$dup.Font.NameFarEast = "Georgia"
$dup.Font.NameBi = "Georgia"
$dup.Font.Name = "Georgia"
$dup.Font.Italic = $true
$dup.Font.Color = 6710886
$dup.Font.Size = 24
$dup.Font.SizeBi = 24

# Give it the same display name as the original so the stylesheet ends up
# with two "Subtitle" entries, matching the duplicated block in the diff.
$dup.NameLocal = "Subtitle"

Write-Output ("Styles.Count=" + $d.Styles.Count)
